$d = $word.ActiveDocument

# Locate the full "<id>p029r_1</id>" text (currently spread across 3 runs
# with different formatting) so we can compute exact sub-range boundaries.
$full = $d.Content
$full.Find.ClearFormatting()
$full.Find.Execute("<id>p029r_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$fullStart = $full.Start
$fullEnd = $full.End

$tag1Len = "<id>".Length
$midLen = "p029r_1".Length

# Range for the first run ("<id>") – this run already carries the
# Courier-New / 7f6000 / sz18 formatting that should apply to the whole tag.
$keepRange = $d.Range($fullStart, $fullStart + $tag1Len)

# Range covering the remaining text ("p029r_1</id>") that currently lives
# in two separately-formatted runs; remove it…
$restRange = $d.Range($fullStart + $tag1Len, $fullEnd)
$restRange.Delete()

# …then re-insert it right after the first run so it inherits that run's
# character formatting and merges into a single run.
$keepRange.InsertAfter("p029r_1</id>")
